$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared strings must be introduced in this exact order so the
# rebuilt sharedStrings table matches the target layout:
# ID(reused), 07/04/2021, 07/04/2022, RGA004, ABCD0RRGA004, ZXC0987RGA004,
# RGA005, ABCD0RRGA005, ZXC0987RGA005
$ws.Range("L2").Value = "'07/04/2021"
$ws.Range("L3").Value = "'07/04/2022"
$ws.Range("U2").Value = "RGA004"
$ws.Range("V2").Value = "ABCD0RRGA004"
$ws.Range("W2").Value = "ZXC0987RGA004"
$ws.Range("U3").Value = "RGA005"
$ws.Range("V3").Value = "ABCD0RRGA005"
$ws.Range("W3").Value = "ZXC0987RGA005"

# Row 2
$ws.Range("F2").Value = 2240451788
$ws.Range("Z2").Value = "No"

# Row 3
$ws.Range("F3").Value = 2240451788

# Selection moves to Z3
$ws.Range("Z3").Select()
